$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.802.12"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "'2.833.35"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'351.24"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("E6").Value = "  +5.71%  "
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.622"
$ws.Range("E9").Value = "  +6.86%  "
$ws.Range("D10").Value = "'40.32"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("D13").Value = "'20.08"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "'7.80"
$ws.Range("E14").Value = "  +3.90%  "
$ws.Range("D15").Value = "'3.271.29"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "'0.975"
$ws.Range("E16").Value = "  +5.70%  "
$ws.Range("D17").Value = "'2.821.12"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "'51.847.07"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "'3.44"
$ws.Range("E19").Value = "  +11.55%  "
$ws.Range("D20").Value = "'7.64"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'13.45"
$ws.Range("D22").Value = "'0.0₃0976"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("D23").Value = "'70.63"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").Value = "'269.74"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("D25").Value = "'2.77"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").Value = "'26.36"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'39.00"
$ws.Range("E29").Value = "  +9.09%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'10.59"
$ws.Range("E30").Value = "  +4.18%  "
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.20"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "'52.79"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "'0.0455"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0901"
$ws.Range("E35").Value = "  +9.07%  "
$ws.Range("D36").Value = "'5.63"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'19.24"
$ws.Range("E38").Value = "  +6.18%  "
$ws.Range("D39").Value = "'3.23"
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").Value = "'121.76"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'22.31"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'2.23"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("E46").Value = "  +8.89%  "
$ws.Range("E47").Value = "  +8.48%  "
$ws.Range("D48").Value = "'2.166.72"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("D49").Value = "'0.997"
$ws.Range("E49").Value = "  +10.40%  "
$ws.Range("D50").Value = "'0.233"
$ws.Range("E50").Value = "  +22.39%  "
$ws.Range("D51").Value = "'0.0323"
$ws.Range("E51").Value = "  +16.72%  "
